# Applies the cryptos-list refresh described in the commit: updated Price (D)
# and Volume(1h) (E) columns across rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4 carries the workbooks plain/default cell style (no explicit format).
# We use it as a reference to restore a D-column cells style after
# temporarily forcing NumberFormat to Text ("@") -- this keeps numeric-looking
# prices (e.g. "540.10", "0.995") stored as literal text instead of being
# auto-coerced into actual numbers by Excel, while leaving the cells visible
# formatting/style exactly as it was before the edit.
$plainStyle = $ws.Range("D4").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.119.82"
$ws.Range("D2").Style = $plainStyle
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.533.75"
$ws.Range("D3").Style = $plainStyle
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.10"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.03"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.563.74"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.53"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  -3.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.363"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.984.45"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.13"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.059.97"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.547.32"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.34"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.07"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.92"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.51"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.440"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.167"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = "  +2.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.993"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.01"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.08"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0800"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.81"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.22"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = "  -6.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.50"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  +5.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "162.12"
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.79"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.51"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.65"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.72"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  -3.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.27"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "303.97"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  -4.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.840"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.73"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.992"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.607"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.84"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.08"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0939"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.27"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0522"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("E51").Value = "  -0.35%  "
